# Fixed bug in main game: the minigame coin-limit threshold changed
# from 30 to 10, in both the English and Vietnamese language rows
# (the "MINIGAME-LIMIT" row of the languages sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 32 holds the "MINIGAME-LIMIT" language row (English in column B,
# Vietnamese in column C). Update the coin threshold from 30 to 10.
$ws.Range("B32").Value = "YOU CAN ONLY PLAY WHEN YOUR COIN IS UNDER 10!"
$ws.Range("C32").Value = "BẠN CHỈ ĐƯỢC CHƠI KHI CÓ ÍT HƠN 10 COIN!"

# Leave the sheet with the same cell selected as the author had when
# they saved the fix.
$ws.Activate()
$ws.Range("C33").Select()
